$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data values
$ws.Range("B2").Value = 331
$ws.Range("B3").Value = 293
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 91

# Delete row 5 (former A5=2,B5=126 row) - shift cells up
$ws.Range("A5:B5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
